# Insert two new price records for "Ají" (Femacal de La Calera) right above
# the existing row 439, shifting the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 439:464 down to 441:466 and make room for the two new rows.
$ws.Range("A439:A440").EntireRow.Insert()

# New row 439
$ws.Cells.Item(439, 1).Value = 3
$ws.Cells.Item(439, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(439, 3).Value = "Coquimbo"
$ws.Cells.Item(439, 4).Value = 44706
$ws.Cells.Item(439, 4).NumberFormat = $ws.Cells.Item(441, 4).NumberFormat
$ws.Cells.Item(439, 5).Value = 5
$ws.Cells.Item(439, 6).Value = 100112021
$ws.Cells.Item(439, 7).Value = "Ají"
$ws.Cells.Item(439, 8).Value = "Americana (o)"
$ws.Cells.Item(439, 9).Value = "Primera"
$ws.Cells.Item(439, 10).Value = 73
$ws.Cells.Item(439, 11).Value = 24000
$ws.Cells.Item(439, 12).Value = 25000
$ws.Cells.Item(439, 13).Value = 24479
$ws.Cells.Item(439, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(439, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(439, 16).Value = 1632
$ws.Cells.Item(439, 17).Value = 15
$ws.Cells.Item(439, 18).Value = "Hortaliza"

# New row 440
$ws.Cells.Item(440, 1).Value = 3
$ws.Cells.Item(440, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(440, 3).Value = "Coquimbo"
$ws.Cells.Item(440, 4).Value = 44706
$ws.Cells.Item(440, 4).NumberFormat = $ws.Cells.Item(441, 4).NumberFormat
$ws.Cells.Item(440, 5).Value = 5
$ws.Cells.Item(440, 6).Value = 100112021
$ws.Cells.Item(440, 7).Value = "Ají"
$ws.Cells.Item(440, 8).Value = "Inferno"
$ws.Cells.Item(440, 9).Value = "Primera"
$ws.Cells.Item(440, 10).Value = 76
$ws.Cells.Item(440, 11).Value = 25000
$ws.Cells.Item(440, 12).Value = 26000
$ws.Cells.Item(440, 13).Value = 25500
$ws.Cells.Item(440, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(440, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(440, 16).Value = 1020
$ws.Cells.Item(440, 17).Value = 25
$ws.Cells.Item(440, 18).Value = "Hortaliza"
